$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New literature records appended 2026-02-03 (GaN CMOS update)
$records = @(
    @('', 'Low Power Design of CMOS Operational Amplifiers for IoT Edge Devices', 2026, 'EWA Publishing', 'Applied and Computational Engineering', 'Zhou, Huayi', '', '10.54254/2755-2721/2026.bj31634', 'https://doi.org/10.54254/2755-2721/2026.bj31634', 'Journal', 'Co-integration', 'Experiment', 'Contacts', '', '', '', 'Low Power Design of CMOS Operational Amplifiers for IoT Edge Devices', 'High', '2026-02-03', ''),
    @('', 'MELEGROS: Monolithic Elephant‐Inspired Gripper with Optical Sensors', 2026, 'Wiley', 'Advanced Science', 'Trunin, Petr; Cafiso, Diana; Nardin, Anderson Brazil; Exley, Trevor; Beccai, Lucia', '', '10.1002/advs.202518878', 'https://doi.org/10.1002/advs.202518878', 'Journal', 'Co-integration', 'Experiment', 'Contacts', '', '', '', 'MELEGROS: Monolithic Elephant‐Inspired Gripper with Optical Sensors', 'High', '2026-02-03', ''),
    @('', 'Pin-Plane Electrical Discharge Driven by a MOSFET DC Current Source', 2026, 'MDPI AG', 'Plasma', 'Perry, Myles; Holoman, Sidmar; Wozniak, Daniel; Dhali, Shirshak Kumar', '', '10.3390/plasma9010005', 'https://doi.org/10.3390/plasma9010005', 'Journal', 'n-FET', 'Experiment', 'Contacts', '', '', '', 'Pin-Plane Electrical Discharge Driven by a MOSFET DC Current Source', 'High', '2026-02-03', ''),
    @('', 'Pin-Plane Electrical Discharge Driven by a MOSFET DC Current Source', 2026, 'MDPI AG', 'Plasma', 'Perry, Myles; Holoman, Sidmar; Wozniak, Daniel; Dhali, Shirshak Kumar', '', '10.3390/plasma9010005', 'https://doi.org/10.3390/plasma9010005', 'Journal', 'n-FET', 'Experiment', 'Contacts', '', '', '', 'Pin-Plane Electrical Discharge Driven by a MOSFET DC Current Source', 'High', '2026-02-03', ''),
    @('', 'Future VLSI Architectures for Neuromorphic Computing, Edge AI and Sustainable Systems', 2027, 'Chandigarh Philosophers', 'International Journal for Multidimensional Research Perspectives', 'R Bhaskar Nihal Varma', '', '10.61877/ijmrp.v3i7.293', 'https://doi.org/10.61877/ijmrp.v3i7.293', 'Journal', 'Co-integration', 'Experiment', 'Contacts', '', '', '', 'Future VLSI Architectures for Neuromorphic Computing, Edge AI and Sustainable Systems', 'High', '2026-02-03', ''),
    @('', 'Low Power Design of CMOS Operational Amplifiers for IoT Edge Devices', 2026, 'EWA Publishing', 'Applied and Computational Engineering', 'Zhou, Huayi', '', '10.54254/2755-2721/2026.bj31634', 'https://doi.org/10.54254/2755-2721/2026.bj31634', 'Journal', 'Co-integration', 'Experiment', 'Contacts', '', '', '', 'Low Power Design of CMOS Operational Amplifiers for IoT Edge Devices', 'High', '2026-02-03', ''),
    @('', 'Interface dipole modulation for gate dielectrics in Field-Effect transistors: a review', 2026, 'Springer Science and Business Media LLC', 'Journal of the Korean Ceramic Society', 'Lim, Wangseop; Kim, Hyojung; Jang, Ho Won', '', '10.1007/s43207-026-00587-5', 'https://doi.org/10.1007/s43207-026-00587-5', 'Journal', 'Co-integration', 'Experiment', 'Contacts', '', '', '', 'Interface dipole modulation for gate dielectrics in Field-Effect transistors: a review', 'High', '2026-02-03', ''),
)

$startRow = 121
for ($i = 0; $i -lt $records.Count; $i++) {
    $r = $startRow + $i
    $rec = $records[$i]
    $ws.Cells.Item($r, 2).Value  = $rec[1]   # Title
    $ws.Cells.Item($r, 3).Value  = $rec[2]   # Year (numeric)
    $ws.Cells.Item($r, 4).Value  = $rec[3]   # Publisher
    $ws.Cells.Item($r, 5).Value  = $rec[4]   # Venue
    $ws.Cells.Item($r, 6).Value  = $rec[5]   # Authors
    $ws.Cells.Item($r, 8).Value  = $rec[7]   # DOI
    $ws.Cells.Item($r, 9).Value  = $rec[8]   # URL
    $ws.Cells.Item($r, 10).Value = $rec[9]   # DocType
    $ws.Cells.Item($r, 11).Value = $rec[10]  # DeviceType
    $ws.Cells.Item($r, 12).Value = $rec[11]  # Method
    $ws.Cells.Item($r, 13).Value = $rec[12]  # EnablerCategory
    $ws.Cells.Item($r, 17).Value = $rec[16]  # EvidenceSnippet
    $ws.Cells.Item($r, 18).Value = $rec[17]  # TagConfidence

    # AddedDate is stored as literal text ("YYYY-MM-DD"), not an Excel date serial,
    # matching the rest of column S. Force text format so Excel does not
    # auto-convert the string into a date value, then restore the default
    # (unstyled) cell style so no extra style index is introduced.
    $cellS = $ws.Cells.Item($r, 19)
    $cellS.NumberFormat = "@"
    $cellS.Value = $rec[18]
    $cellS.Style = "Normal"
}

